$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newUrl = "https://web.archive.org/web/19970719105949/http://dewey.rug.ac.be/YFCF/HotSonic.html"

# B33 currently holds the (dead) archived HotSonic URL; replace it with the
# newly-found archive.org snapshot and turn it into a real hyperlink.
$cell = $ws.Range("B33")
$cell.Value = $newUrl

$ws.Hyperlinks.Add($cell, $newUrl)

# Reflect the cursor position left behind in the saved file.
$ws.Range("D15").Select()

$wb.Save()
